$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G4").Value = "mariam.noureldin@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G5").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G8").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G15").Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G17").Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G18").Value = "marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G19").Value = "rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G21").Value = "mariam.noureldin@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G22").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G25").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G32").Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G34").Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G37").Value = "System, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G45").Value = "System, backup@backdoor.com, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G51").Value = "naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G52").Value = "abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("G54").Value = "System, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G62").Value = "System, backup@backdoor.com, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G68").Value = "naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G69").Value = "abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("G70").Value = "rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G72").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G73").Value = "mariam.noureldin@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G76").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G77").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G83").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G85").Value = "Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G86").Value = "marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G87").Value = "rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G89").Value = "mariam.noureldin@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G92").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G98").Value = "Walaa.h.ghanima@med.asu.edu.eg, nourhanhosni@med.asu.edu.eg, user@user.com, nourhanmohamed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G100").Value = "Aya.alaa-eldein@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G102").Value = "Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G105").Value = "System, backup@backdoor.com"
$ws.Range("G106").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G110").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G116").Value = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G117").Value = "enas.omran@med.asu.edu.eg, takyas@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg`n, user@user.com"
$ws.Range("G119").Value = "marina_atef@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G122").Value = "System, backup@backdoor.com"
$ws.Range("G123").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G127").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G133").Value = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G134").Value = "enas.omran@med.asu.edu.eg, takyas@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg`n, user@user.com"
$ws.Range("G136").Value = "marina_atef@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G144").Value = "Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Range("G150").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G151").Value = "System, randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
